$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (cell styles) from an existing data row down to the
# new row 38 before we overwrite any values.
$ws.Range("A33:B33").Copy()
$ws.Range("A38:B38").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Replace the hypercholesterolemia rows (33-37) with heart-failure rows,
# and populate the new Multiple Sclerosis row (38).

$ws.Cells.Item(33,1).Value = "I50"
$ws.Cells.Item(33,2).Value = "hf"
$ws.Cells.Item(33,3).Value = "heart failure"

$ws.Cells.Item(34,1).Value = "I110"
$ws.Cells.Item(34,2).Value = "hf"
$ws.Cells.Item(34,3).Value = "HTN heart disease with CHF"

$ws.Cells.Item(35,1).Value = "I130"
$ws.Cells.Item(35,2).Value = "hf"
$ws.Cells.Item(35,3).Value = "HTN heart and renal disease with CHF"

$ws.Cells.Item(36,1).Value = "I132"
$ws.Cells.Item(36,2).Value = "hf"
$ws.Cells.Item(36,3).Value = "HTN with heart and renal disease with both CKD and CHF"

$ws.Cells.Item(37,1).Value = "I42"
$ws.Cells.Item(37,2).Value = "hf"
$ws.Cells.Item(37,3).Value = "Cardiomyopathies"

$ws.Cells.Item(38,1).Value = "G35"
$ws.Cells.Item(38,2).Value = "ms"
$ws.Cells.Item(38,3).Value = "Multiple Sclerosis"

$ws.Range("C39").Select()
